# Apply the "Office Theme" design to the deck's (single) slide master,
# and reset the summary table on slide 16 back to the default table
# style that comes with a freshly-applied theme.
#
# The presentation currently uses the "Integral" theme on its slide
# master (ppt/theme/theme1.xml). Re-coloring it to the stock "Office
# Theme" palette (font scheme + format scheme are already identical
# between the two themes, so only the 10 non-black/white scheme colors
# actually change) reproduces the effective result of picking the
# "Office Theme" design from the Design gallery.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$clrs = $master.Theme.ThemeColorScheme

# MsoThemeColorSchemeIndex order: Dark1, Light1, Dark2, Light2,
# Accent1..Accent6, Hyperlink, FollowedHyperlink.
# .RGB takes/returns a VBA-style 0x00BBGGRR packed long.
$clrs.Item(1).RGB  = 0          # dk1      #000000 (unchanged)
$clrs.Item(2).RGB  = 16777215   # lt1      #FFFFFF (unchanged)
$clrs.Item(3).RGB  = 6968388    # dk2      #44546A
$clrs.Item(4).RGB  = 15132391   # lt2      #E7E6E6
$clrs.Item(5).RGB  = 13998939   # accent1  #5B9BD5
$clrs.Item(6).RGB  = 3243501    # accent2  #ED7D31
$clrs.Item(7).RGB  = 10855845   # accent3  #A5A5A5
$clrs.Item(8).RGB  = 49407      # accent4  #FFC000
$clrs.Item(9).RGB  = 12874308   # accent5  #4472C4
$clrs.Item(10).RGB = 4697456    # accent6  #70AD47
$clrs.Item(11).RGB = 12673797   # hlink    #0563C1
$clrs.Item(12).RGB = 7491477    # folHlink #954F72

# The theme switch resets the table on slide 16 (the cash-flow recap
# table) to the new theme's default table style ("Themed Style 2 -
# Accent 1") instead of the old custom "Table_0" style.
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{9810380A-5DFF-4B4F-9CFA-EAD38923DC8B}")
    }
}
